$d = $word.ActiveDocument

# Locate the paragraph that ends with "Terminates abnormally." (capital T),
# the last bullet item under "Finally block - Usage".
$rng = $d.Content
$found = $rng.Find.Execute("Terminates abnormally.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph 'Terminates abnormally.'"
}
$para = $rng.Paragraphs(1)

# Insert a brand new paragraph right after it.
$para.Range.InsertParagraphAfter()
$newPara = $para.Next()

# The source paragraph is a numbered/bulleted ListParagraph; the new
# paragraph in the diff has no pPr/list formatting at all, so strip
# whatever got inherited.
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = "Normal"

# Collapse to the (now plain, empty) new paragraph's start and inject the
# exact run structure -- including bold runs and proofErr spell/grammar
# markers -- via raw WordprocessingML.
$target = $newPara.Range
$target.Collapse(1)

$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>User defined exception</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>&#8211; java allows us to</w:t></w:r><w:r><w:t xml:space="preserve"> create our own exception</w:t></w:r><w:r><w:t xml:space="preserve">, which is essentially a derived class of exception. To create our own exception, we must first create a class that extends the exception class and represents user defined </w:t></w:r><w:r><w:t xml:space="preserve">exceptions. We must pass the string to constructor if the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>superclass</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, which is obtained by calling the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getMessage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) function if the newly created object. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xmlFrag) | Out-Null

Write-Output "Inserted 'User defined exception' paragraph."
